$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header id row (row 1)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Update CON data row (row 2)
$ws.Range("B2").Value = 12.637594603067761
$ws.Range("C2").Value = 11.088586800610617
$ws.Range("D2").Value = 12.235817162159993
$ws.Range("E2").Value = 12.009443207668815

# Update STR data row (row 3)
$ws.Range("B3").Value = 12.502730745004795
$ws.Range("C3").Value = 10.319760442308223
$ws.Range("D3").Value = 12.957339580952244
$ws.Range("E3").Value = 10.672956785928051

# Update selection to match new active range
$ws.Range("B1:E3").Select()
